$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "Den of Thieves"
$ws.Range("C31").Value = "James B. Stewart"
$ws.Range("D31").Value = "A #1 bestseller from coast to coast, Den of Thieves tells the full story of the insider-trading scandal that nearly destroyed Wall Street, the men who pulled it off, and the chase that finally brought them to justice."
$ws.Range("E31").Value = "30 Den of Thieves.jpg"
$ws.Range("F31").Value = 45560
